$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "bus" sheet: update the x/y graphical-position columns (I, J) for the
#    three buses. These values are stored as TEXT in the workbook (not
#    numbers), so a leading apostrophe is used to force Excel to keep the
#    value as text; the cell style is then reset to "Normal" so no stray
#    quote-prefix / number-format styling is left behind on the cell.
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("bus")

# Bus 1 (row 2): x -773.0 -> -1191.0, y -466.0 -> -716.0
$wsBus.Range("I2").Value = "'-1191.0"
$wsBus.Range("I2").Style = "Normal"
$wsBus.Range("J2").Value = "'-716.0"
$wsBus.Range("J2").Style = "Normal"

# Bus 2 (row 3): x -1039.0 -> -1457.0, y -611.0 -> -861.0
$wsBus.Range("I3").Value = "'-1457.0"
$wsBus.Range("I3").Style = "Normal"
$wsBus.Range("J3").Value = "'-861.0"
$wsBus.Range("J3").Style = "Normal"

# Bus 3 (row 4): x -901.0 -> -1319.0, y -324.0 -> -574.0
$wsBus.Range("I4").Value = "'-1319.0"
$wsBus.Range("I4").Style = "Normal"
$wsBus.Range("J4").Value = "'-574.0"
$wsBus.Range("J4").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) "battery" sheet: rename the battery from "batt" to "batt1@Bus 3"
# ---------------------------------------------------------------------------
$wsBattery = $wb.Worksheets.Item("battery")
$wsBattery.Range("B2").Value = "batt1@Bus 3"

# ---------------------------------------------------------------------------
# 3) "branch" sheet: fix DCOPF bug to consider storage by giving the
#    branches a proper reactance (X) - and for branch row 4 also resistance
#    (R) - instead of the ~0 placeholder value (9.999999999999999e-21).
# ---------------------------------------------------------------------------
$wsBranch = $wb.Worksheets.Item("branch")

$wsBranch.Range("J2").Value = 0.05
$wsBranch.Range("J3").Value = 0.08
$wsBranch.Range("I4").Value = 0.01
$wsBranch.Range("J4").Value = 0.06
